$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F3").Value = 420
$ws.Range("F4").Value = 154
$ws.Range("F6").Value = 3845
$ws.Range("F9").Value = 76
$ws.Range("F10").Value = 3114
$ws.Range("F13").Value = 2301
$ws.Range("G13").Value = 55
$ws.Range("F16").Value = 89
$ws.Range("F18").Value = 4
$ws.Range("F20").Value = 206
$ws.Range("F23").Value = 376
$ws.Range("F24").Value = 657
$ws.Range("F27").Value = 5
$ws.Range("F28").Value = 1295
$ws.Range("F29").Value = 128
$ws.Range("F31").Value = 11
$ws.Range("F33").Value = 48
$ws.Range("F34").Value = 4277
$ws.Range("F35").Value = 3985
$ws.Range("F39").Value = 1122
$ws.Range("F41").Value = 466
$ws.Range("F42").Value = 11
$ws.Range("F44").Value = 163
$ws.Range("F45").Value = 128

$ws = $wb.Worksheets.Item(2)
$ws.Range("F7").Value = 24
$ws.Range("F15").Value = 207

$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 145
$ws.Range("F4").Value = 2285

$ws = $wb.Worksheets.Item(4)
$ws.Range("F3").Value = 145
$ws.Range("F4").Value = 420
$ws.Range("F7").Value = 154
$ws.Range("F9").Value = 3845
$ws.Range("F12").Value = 76
$ws.Range("F13").Value = 3114
$ws.Range("F15").Value = 2301
$ws.Range("G15").Value = 55
$ws.Range("F17").Value = 89
$ws.Range("F18").Value = 4
$ws.Range("F22").Value = 376
$ws.Range("F23").Value = 657
$ws.Range("F26").Value = 1295
$ws.Range("F27").Value = 128
$ws.Range("F30").Value = 48
$ws.Range("F31").Value = 24
$ws.Range("F32").Value = 4277
$ws.Range("F38").Value = 466
$ws.Range("F40").Value = 11
$ws.Range("F44").Value = 163
$ws.Range("F49").Value = 207
